$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 265.57144
$ws.Range("I8").Value = 276.5
$ws.Range("K8").Value = 829.5
$ws.Range("M8").Value = -690.5
$ws.Range("H74").Value = 5920.5713
$ws.Range("I74").Value = 5662.6665
$ws.Range("J74").Value = 5990.909
$ws.Range("K74").Value = 5662.6665
$ws.Range("L74").Value = 5990.909
$ws.Range("M74").Value = -4726.6665
$ws.Range("N74").Value = -7862.909
$ws.Range("H77").Value = 5920.5713
$ws.Range("I77").Value = 5662.6665
$ws.Range("J77").Value = 5990.909
$ws.Range("K77").Value = 28313.3325
$ws.Range("L77").Value = 29954.545
$ws.Range("M77").Value = -23633.3325
$ws.Range("N77").Value = -39314.545
$ws.Range("H92").Value = 659.23334
$ws.Range("I92").Value = 635.9231
$ws.Range("K92").Value = 635.9231
$ws.Range("M92").Value = 612.0769
$ws.Range("H101").Value = 264.33334
$ws.Range("I101").Value = 264.33334
$ws.Range("K101").Value = 793.0000200000001
$ws.Range("M101").Value = 828.9999799999999
$ws.Range("H137").Value = 3861.6875
$ws.Range("I137").Value = 1377.0741
$ws.Range("J137").Value = 17278.6
$ws.Range("K137").Value = 4131.2223
$ws.Range("L137").Value = 51835.8
$ws.Range("M137").Value = -1581.2223
$ws.Range("N137").Value = -56935.8
$ws.Range("H138").Value = 6800.8296
$ws.Range("J138").Value = 10156.241
$ws.Range("L138").Value = 30468.723
$ws.Range("N138").Value = -40748.723

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 4645.1665
$ws.Range("I2").Value = 3716.625
$ws.Range("J2").Value = 6502.25
$ws.Range("K2").Value = 3716.625
$ws.Range("L2").Value = 6502.25
$ws.Range("M2").Value = -3603.625
$ws.Range("N2").Value = -6728.25
$ws.Range("H32").Value = 7764.026
$ws.Range("I32").Value = 6653.7334
$ws.Range("J32").Value = 11465
$ws.Range("K32").Value = 6653.7334
$ws.Range("L32").Value = 11465
$ws.Range("M32").Value = -6366.7334
$ws.Range("N32").Value = -12039
$ws.Range("H55").Value = 40715.25
$ws.Range("J55").Value = 47620.332
$ws.Range("L55").Value = 47620.332
$ws.Range("N55").Value = -48250.332
$ws.Range("H74").Value = 22224818
$ws.Range("I74").Value = 55557570
$ws.Range("J74").Value = 2985.111
$ws.Range("K74").Value = 55557570
$ws.Range("L74").Value = 2985.111
$ws.Range("M74").Value = -55556696
$ws.Range("N74").Value = -4733.111
$ws.Range("H77").Value = 22224818
$ws.Range("I77").Value = 55557570
$ws.Range("J77").Value = 2985.111
$ws.Range("K77").Value = 277787850
$ws.Range("L77").Value = 14925.555
$ws.Range("M77").Value = -277783482
$ws.Range("N77").Value = -23661.555
$ws.Range("H97").Value = 2022.826
$ws.Range("I97").Value = 1764.7894
$ws.Range("K97").Value = 1764.7894
$ws.Range("M97").Value = -1268.7894
$ws.Range("H116").Value = 4645.1665
$ws.Range("I116").Value = 3716.625
$ws.Range("J116").Value = 6502.25
$ws.Range("K116").Value = 3716.625
$ws.Range("L116").Value = 6502.25
$ws.Range("M116").Value = -1422.625
$ws.Range("N116").Value = -11090.25
$ws.Range("H132").Value = 21308366
$ws.Range("I132").Value = 2011.9667
$ws.Range("K132").Value = 6035.9001
$ws.Range("M132").Value = -3505.9001

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 4645.1665
$ws.Range("I3").Value = 3716.625
$ws.Range("J3").Value = 6502.25
$ws.Range("K3").Value = 3716.625
$ws.Range("L3").Value = 6502.25
$ws.Range("M3").Value = -3602.625
$ws.Range("N3").Value = -6730.25
$ws.Range("H37").Value = 1893.75
$ws.Range("I37").Value = 1087.5
$ws.Range("J37").Value = 2700
$ws.Range("K37").Value = 1087.5
$ws.Range("L37").Value = 2700
$ws.Range("M37").Value = -950.5
$ws.Range("N37").Value = -2974
$ws.Range("H86").Value = 8580.772000000001
$ws.Range("I86").Value = 9809
$ws.Range("K86").Value = 9809
$ws.Range("M86").Value = -8686
$ws.Range("H89").Value = 8580.772000000001
$ws.Range("I89").Value = 9809
$ws.Range("K89").Value = 49045
$ws.Range("M89").Value = -43429
$ws.Range("H94").Value = 1748.4375
$ws.Range("I94").Value = 1521.1538
$ws.Range("J94").Value = 2733.3333
$ws.Range("K94").Value = 1521.1538
$ws.Range("L94").Value = 2733.3333
$ws.Range("M94").Value = -1070.1538
$ws.Range("N94").Value = -3635.3333
$ws.Range("H105").Value = 6632.5806
$ws.Range("I105").Value = 14558.223
$ws.Range("J105").Value = 3390.2727
$ws.Range("K105").Value = 14558.223
$ws.Range("L105").Value = 3390.2727
$ws.Range("M105").Value = -12811.223
$ws.Range("N105").Value = -6884.2727

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 11368940
$ws.Range("J31").Value = 17863202
$ws.Range("L31").Value = 17863202
$ws.Range("N31").Value = -17863792
$ws.Range("H34").Value = 11368940
$ws.Range("J34").Value = 17863202
$ws.Range("L34").Value = 17863202
$ws.Range("N34").Value = -17863606
$ws.Range("H131").Value = 22598.5
$ws.Range("J131").Value = 25798
$ws.Range("L131").Value = 25798
$ws.Range("N131").Value = -35878
$ws.Range("H132").Value = 3181.4211
$ws.Range("I132").Value = 3120.4707
$ws.Range("J132").Value = 3699.5
$ws.Range("K132").Value = 9361.4121
$ws.Range("L132").Value = 11098.5
$ws.Range("M132").Value = -6831.4121
$ws.Range("N132").Value = -16158.5
$ws.Range("H134").Value = 3402.5715
$ws.Range("I134").Value = 2785.0908
$ws.Range("J134").Value = 5666.6665
$ws.Range("K134").Value = 8355.2724
$ws.Range("L134").Value = 16999.9995
$ws.Range("M134").Value = -5820.2724
$ws.Range("N134").Value = -22069.9995
$ws.Range("H141").Value = 130257.836
$ws.Range("J141").Value = 130257.836
$ws.Range("L141").Value = 130257.836
$ws.Range("N141").Value = -140617.836

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H128").Value = 196515
$ws.Range("I128").Value = 196515
$ws.Range("K128").Value = 589545
$ws.Range("M128").Value = -584565
$ws.Range("H131").Value = 1326
$ws.Range("I131").Value = 1038.4445
$ws.Range("J131").Value = 1695.7142
$ws.Range("K131").Value = 3115.3335
$ws.Range("L131").Value = 5087.142599999999
$ws.Range("M131").Value = 1924.6665
$ws.Range("N131").Value = -15167.1426
$ws.Range("H134").Value = 6433
$ws.Range("J134").Value = 11516.5
$ws.Range("L134").Value = 34549.5
$ws.Range("N134").Value = -44689.5

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3694.3333
$ws.Range("I102").Value = 2053
$ws.Range("K102").Value = 2053
$ws.Range("M102").Value = -431
$ws.Range("H107").Value = 947.5
$ws.Range("I107").Value = 947.5
$ws.Range("K107").Value = 947.5
$ws.Range("M107").Value = 972.5
$ws.Range("H122").Value = 19233070
$ws.Range("I122").Value = 2068.8
$ws.Range("J122").Value = 83336410
$ws.Range("K122").Value = 6206.400000000001
$ws.Range("L122").Value = 250009230
$ws.Range("M122").Value = -3756.400000000001
$ws.Range("N122").Value = -250014130
$ws.Range("H123").Value = 41121.145
$ws.Range("J123").Value = 36212
$ws.Range("L123").Value = 36212
$ws.Range("N123").Value = -41112
$ws.Range("H126").Value = 69274.56
$ws.Range("I126").Value = 97035.73
$ws.Range("J126").Value = 8200
$ws.Range("K126").Value = 291107.19
$ws.Range("L126").Value = 24600
$ws.Range("M126").Value = -288637.19
$ws.Range("N126").Value = -29540

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5073.778
$ws.Range("I7").Value = 5023.5713
$ws.Range("J7").Value = 5249.5
$ws.Range("K7").Value = 5023.5713
$ws.Range("L7").Value = 5249.5
$ws.Range("M7").Value = -4911.5713
$ws.Range("N7").Value = -5473.5
$ws.Range("H40").Value = 2670.6667
$ws.Range("I40").Value = 3006.1667
$ws.Range("J40").Value = 1999.6666
$ws.Range("K40").Value = 3006.1667
$ws.Range("L40").Value = 1999.6666
$ws.Range("M40").Value = -2870.1667
$ws.Range("N40").Value = -2271.6666
$ws.Range("H93").Value = 563077.7
$ws.Range("I93").Value = 2119.0527
$ws.Range("K93").Value = 2119.0527
$ws.Range("M93").Value = -871.0527000000002
$ws.Range("H122").Value = 3575595.5
$ws.Range("I122").Value = 3922
$ws.Range("J122").Value = 12504779
$ws.Range("K122").Value = 11766
$ws.Range("L122").Value = 37514337
$ws.Range("M122").Value = -9316
$ws.Range("N122").Value = -37519237
$ws.Range("H126").Value = 5073.778
$ws.Range("I126").Value = 5023.5713
$ws.Range("J126").Value = 5249.5
$ws.Range("K126").Value = 15070.7139
$ws.Range("L126").Value = 15748.5
$ws.Range("M126").Value = -12600.7139
$ws.Range("N126").Value = -20688.5
$ws.Range("H131").Value = 69811.5
$ws.Range("J131").Value = 86975
$ws.Range("L131").Value = 86975
$ws.Range("N131").Value = -97055

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 6405.5386
$ws.Range("I81").Value = 4699.6665
$ws.Range("J81").Value = 6917.3
$ws.Range("K81").Value = 9399.333000000001
$ws.Range("L81").Value = 13834.6
$ws.Range("M81").Value = -8338.333000000001
$ws.Range("N81").Value = -15956.6
$ws.Range("H84").Value = 6405.5386
$ws.Range("I84").Value = 4699.6665
$ws.Range("J84").Value = 6917.3
$ws.Range("K84").Value = 46996.665
$ws.Range("L84").Value = 69173
$ws.Range("M84").Value = -41692.665
$ws.Range("N84").Value = -79781
$ws.Range("H86").Value = 12535624
$ws.Range("J86").Value = 32998
$ws.Range("L86").Value = 32998
$ws.Range("N86").Value = -35244
$ws.Range("H89").Value = 12535624
$ws.Range("J89").Value = 32998
$ws.Range("L89").Value = 164990
$ws.Range("N89").Value = -176222
$ws.Range("H100").Value = 48096630
$ws.Range("I100").Value = 72144450
$ws.Range("K100").Value = 144288900
$ws.Range("M100").Value = -144288359
$ws.Range("H122").Value = 9526947
$ws.Range("I122").Value = 1572.3334
$ws.Range("J122").Value = 33340384
$ws.Range("K122").Value = 4717.0002
$ws.Range("L122").Value = 100021152
$ws.Range("M122").Value = -2267.0002
$ws.Range("N122").Value = -100026052
